# Apply the cryptos-list refresh: updated prices/volumes for all rows,
# plus a rank swap for rows 42-45 (TrustWalletToken/FraxShare and
# RenderToken/PaxDollar trade places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string (e.g. "0.9993") to be stored as
# text, matching the inlineStr cells already in the sheet, instead of
# letting Excel auto-convert it to a real number.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '30.015.89'
$ws.Range("E2").Value = '  -0.32%  '

# Row 3
$ws.Range("D3").Value = '1.899.30'
$ws.Range("E3").Value = '  -0.79%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.9993'
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
Set-TextValue $ws.Range("D5") '0.8367'
$ws.Range("E5").Value = '  +4.74%  '

# Row 6
Set-TextValue $ws.Range("D6") '242.10'
$ws.Range("E6").Value = '  -0.58%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.9994'
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3283'
$ws.Range("E8").Value = '  +2.65%  '

# Row 9
Set-TextValue $ws.Range("D9") '26.55'
$ws.Range("E9").Value = '  +0.76%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.07046'
$ws.Range("E10").Value = '  +1.21%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.08084'
$ws.Range("E11").Value = '  +1.00%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.7635'
$ws.Range("E12").Value = '  +1.43%  '

# Row 13
$ws.Range("D13").Value = '1.906.44'
$ws.Range("E13").Value = '  -0.24%  '

# Row 14
Set-TextValue $ws.Range("D14") '5.260'
$ws.Range("E14").Value = '  +0.36%  '

# Row 15
Set-TextValue $ws.Range("D15") '92.38'
$ws.Range("E15").Value = '  -1.40%  '

# Row 16
$ws.Range("D16").Value = '30.017.59'
$ws.Range("E16").Value = '  -0.38%  '

# Row 17
Set-TextValue $ws.Range("D17") '14.12'
$ws.Range("E17").Value = '  +0.34%  '

# Row 18
Set-TextValue $ws.Range("D18") '5.858'
$ws.Range("E18").Value = '  -2.29%  '

# Row 19
Set-TextValue $ws.Range("D19") '244.14'
$ws.Range("E19").Value = '  -2.09%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.000007765'
$ws.Range("E20").Value = '  -0.82%  '

# Row 21
Set-TextValue $ws.Range("D21") '0.9996'
$ws.Range("E21").Value = '  -0.06%  '

# Row 22
$ws.Range("D22").Value = '2.148.36'
$ws.Range("E22").Value = '  -0.94%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.9992'
$ws.Range("E23").Value = '  -0.06%  '

# Row 24
Set-TextValue $ws.Range("D24") '6.968'
$ws.Range("E24").Value = '  -0.03%  '

# Row 25
Set-TextValue $ws.Range("D25") '0.1745'
$ws.Range("E25").Value = '  +24.31%  '

# Row 26
Set-TextValue $ws.Range("D26") '9.292'
$ws.Range("E26").Value = '  -0.47%  '

# Row 27
Set-TextValue $ws.Range("D27") '165.41'
$ws.Range("E27").Value = '  -2.43%  '

# Row 28
Set-TextValue $ws.Range("D28") '18.94'
$ws.Range("E28").Value = '  -0.34%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.094'
$ws.Range("E29").Value = '  +1.60%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.362'
$ws.Range("E30").Value = '  -2.09%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.514'
$ws.Range("E31").Value = '  -1.00%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.05906'
$ws.Range("E32").Value = '  +8.16%  '

# Row 33
Set-TextValue $ws.Range("D33") '4.285'
$ws.Range("E33").Value = '  -1.71%  '

# Row 34
Set-TextValue $ws.Range("D34") '4.073'
$ws.Range("E34").Value = '  -1.39%  '

# Row 35
Set-TextValue $ws.Range("D35") '1.269'
$ws.Range("E35").Value = '  +0.08%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.7324'
$ws.Range("E36").Value = '  -0.96%  '

# Row 37
Set-TextValue $ws.Range("D37") '2.715'
$ws.Range("E37").Value = '  -0.67%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.01919'
$ws.Range("E38").Value = '  -0.89%  '

# Row 39
Set-TextValue $ws.Range("D39") '2.780'
$ws.Range("E39").Value = '  -0.57%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.4447'
$ws.Range("E40").Value = '  -0.47%  '

# Row 41
Set-TextValue $ws.Range("D41") '72.81'
$ws.Range("E41").Value = '  -0.66%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D42") '0.8579'
$ws.Range("E42").Value = '  +2.77%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D43") '5.871'
$ws.Range("E43").Value = '  -5.10%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D44") '1.908'
$ws.Range("E44").Value = '  +0.06%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D45") '0.9991'
$ws.Range("E45").Value = '  -0.17%  '

# Row 46
$ws.Range("E46").Value = '  +0.91%  '

# Row 47
Set-TextValue $ws.Range("D47") '7.561'
$ws.Range("E47").Value = '  -0.68%  '

# Row 48
Set-TextValue $ws.Range("D48") '9.811'
$ws.Range("E48").Value = '  -0.72%  '

# Row 49
Set-TextValue $ws.Range("D49") '996.64'
$ws.Range("E49").Value = '  +1.77%  '

# Row 50
$ws.Range("D50").Value = '2.047.98'
$ws.Range("E50").Value = '  -0.83%  '

# Row 51
Set-TextValue $ws.Range("D51") '1.516'
$ws.Range("E51").Value = '  +0.33%  '
